# CI Build: Synchronisation of the project's models with the google sheet metadata version 10142
# Add preparationTechnique property to Product (and related subtype/usage sheets)
# Add logo property to FundingSource

$wb = $excel.ActiveWorkbook

# Sheets that gain a new "preparationTechnique" column, inserted immediately
# before the existing "accessPointUrl" column. Column index is 1-based.
$preparationTechniqueTargets = @(
    @{ Sheet = "Product";      Col = 8  },
    @{ Sheet = "Antibody";     Col = 13 },
    @{ Sheet = "Hybridoma";    Col = 14 },
    @{ Sheet = "Protein";      Col = 22 },
    @{ Sheet = "NucleicAcid";  Col = 22 },
    @{ Sheet = "DetectionKit"; Col = 12 },
    @{ Sheet = "Bundle";       Col = 9  },
    @{ Sheet = "Virus";        Col = 29 },
    @{ Sheet = "Bacterium";    Col = 26 },
    @{ Sheet = "Fungus";       Col = 26 },
    @{ Sheet = "Protozoan";    Col = 26 },
    @{ Sheet = "Viroid";       Col = 26 },
    @{ Sheet = "Prion";        Col = 26 }
)

foreach ($target in $preparationTechniqueTargets) {
    $ws = $wb.Worksheets.Item($target.Sheet)
    $ws.Columns.Item($target.Col).Insert()
    $cell = $ws.Cells.Item(1, $target.Col)
    $cell.Value = "preparationTechnique"
}

# FundingSource gains a new "logo" column, inserted immediately before the
# existing "keyword" column.
$wsFunding = $wb.Worksheets.Item("FundingSource")
$wsFunding.Columns.Item(9).Insert()
$wsFunding.Cells.Item(1, 9).Value = "logo"
